# Update existing rows 2-45 with new B/C/D values, and append new rows 46-69
# (A column = row-2 index, B/C/D = data values), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.645206191993366, 8.858668200845866, 0.4440932171299949),
    @(3, 4.383809600599529, 1.382641365688188, 0.3289346777634581),
    @(4, 4.847638242716942, 5.888908002725834, 0.7054606029343404),
    @(5, 6.37223141104835, 7.117657876314829, 0.5620629427854495),
    @(6, 6.766458026411556, 1.248686260683577, 0.3712522531804548),
    @(7, 12.49679203521436, 2.608616447873291, 0.3171336951303816),
    @(8, 15.18132571376524, 7.803280427098616, 0.3482744938681708),
    @(9, 16.0965490438033, 4.5344284782238, 0.4656379717125014),
    @(10, 17.39387357945646, 4.073331090778641, 0.1826336386989122),
    @(11, 20.51691641544009, 2.638912571447953, 0.5030849136238591),
    @(12, 21.25921810457296, 6.752993861897306, 0.7131511897723457),
    @(13, 26.32249086448913, 6.841613101383336, 0.5990683665144856),
    @(14, 26.87396203536685, 2.164824338768115, 0.1695915333235668),
    @(15, 32.93449667377096, 2.517791272377724, 0.1996924488484113),
    @(16, 34.0272698158478, 1.512464596366264, 0.3257071581399337),
    @(17, 35.05815693838623, 5.468755042558919, 0.3459356317888655),
    @(18, 36.53744271072388, 6.424061610919515, 0.5660353622638917),
    @(19, 37.10028371783362, 5.892855347159736, 0.3370742401896643),
    @(20, 37.94609828741224, 3.128345679075647, 0.3334217544812892),
    @(21, 38.06010790228134, 3.589218317300853, 0.4593518133396342),
    @(22, 38.40401611038392, 8.847637558557913, 0.4947951705066551),
    @(23, 38.64532519864624, 1.121299289936611, 0.3266422111505823),
    @(24, 38.75533702195085, 6.958790388587532, 0.7038443241530404),
    @(25, 39.42329033841006, 5.081219019909729, 0.7927182016156993),
    @(26, 39.53923724316722, 2.315057034167105, 0.3973057821477242),
    @(27, 39.79631869285787, 9.349941952163016, 0.6021271840216795),
    @(28, 42.1521477778083, 3.738655033172221, 0.7922226568226339),
    @(29, 45.24739522238139, 4.940961599263082, 0.5763937101637046),
    @(30, 45.70923068005656, 3.95843204787689, 0.2535721347299477),
    @(31, 46.28325937043653, 5.769691346838951, 0.450081875232447),
    @(32, 48.31820076269299, 8.945481260757317, 0.5195370194071269),
    @(33, 48.44556645546279, 2.147899360231266, 0.428300133800611),
    @(34, 48.52343564834958, 4.167361037283245, 0.4829081202743919),
    @(35, 48.82553677705712, 3.576305893893265, 0.5882196793411353),
    @(36, 50.30248362632528, 3.337955470427211, 0.3216242038921238),
    @(37, 50.86219542283866, 3.956664248563747, 0.6359673210002313),
    @(38, 53.1641258154374, 8.243109264058571, 0.4213770678295656),
    @(39, 53.88793558981771, 4.63105475207677, 0.4572020246092904),
    @(40, 54.54826684180106, 1.356010027097246, 0.2639647832122718),
    @(41, 56.27652665119858, 5.513855829400645, 0.6434772185748003),
    @(42, 59.3033849387793, 1.844219565533461, 0.2454813566728519),
    @(43, 59.82904322343599, 7.735191653280012, 0.620321975776104),
    @(44, 60.16084019636168, 3.522244768610373, 0.5272825998879047),
    @(45, 61.71501396626988, 5.947976784914683, 0.4412564875950641),
    @(46, 62.51355431783183, 1.31825628611009, 0.2867965375498073),
    @(47, 64.76271612924744, 2.106664488603526, 0.2440283094882096),
    @(48, 65.71118310506517, 4.939281618008239, 0.6908720199816579),
    @(49, 67.04118931534715, 4.313215490149245, 0.5098213271858253),
    @(50, 71.10966172417534, 6.238241857397202, 0.6137776012408273),
    @(51, 72.49902512663589, 4.264546063320592, 0.4069104197161962),
    @(52, 72.84200276627401, 4.508505855531086, 0.4007517000747983),
    @(53, 73.6226621543851, 6.081822720436659, 0.6374089664423284),
    @(54, 73.77023551779313, 6.880731461154667, 0.5074374292925939),
    @(55, 74.3955741971406, 2.852469498987257, 0.1924847904321531),
    @(56, 74.48370499181156, 4.929428217862749, 0.7483200854395077),
    @(57, 77.87296389163599, 2.844719514106053, 0.344646813939186),
    @(58, 79.24974140182462, 3.073408159722917, 0.3883076727536454),
    @(59, 79.55495946545591, 2.720411975841293, 0.350417183644694),
    @(60, 81.80872007210753, 6.127069127655782, 0.6429914093265965),
    @(61, 83.89014339941453, 5.451138367438002, 0.7472319108590775),
    @(62, 84.01678091554064, 4.873240246308664, 0.6061867032193685),
    @(63, 85.99683130084721, 3.985522459530099, 0.6610669502217461),
    @(64, 86.75477267356237, 6.07795040056809, 0.5619106822635931),
    @(65, 89.55641069374384, 3.025063499180566, 0.3129615517519923),
    @(66, 90.59107136184295, 7.567334824958841, 0.6160474318631286),
    @(67, 96.20586597829954, 5.84155445908424, 0.5339030016565086),
    @(68, 98.33242109002627, 7.393803162446424, 0.5199333714701926),
    @(69, 98.34569998117593, 2.051873069320413, 0.2803364471094464)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $b = $entry[1]
    $c = $entry[2]
    $d = $entry[3]

    if ($r -gt 45) {
        # New row: set the index in column A and copy column-A formatting
        # (bold/bordered/centered style) from the row above.
        $ws.Cells.Item($r, 1).Value = ($r - 2)
    }

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}

# Copy the formatting of A45 (bold, bordered, centered) down through the
# newly-added A46:A69 cells so they match the existing index column style.
$styleSource = $ws.Range("A45")
$styleSource.Copy()
$styleTarget = $ws.Range("A46:A69")
$styleTarget.PasteSpecial(-4122)

Write-Output "done"
